# [MOSIP-14369] Fix: boolean values
#
# The "is_global_working" (column E) and "is_active" (column F) columns were
# populated with =TRUE()/=FALSE() boolean formulas. They should instead just
# hold the literal text "TRUE"/"FALSE", matching the rest of the masterdata
# upload sheets (which store booleans as plain text so the importer doesn't
# choke on real Excel boolean cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 and row 9 happen to pick up a dedicated "quoted text" cell style in
# the original edit (because those two FALSE cells were re-typed by hand),
# so they are handled separately below.
$quotedStyleRows = @(8, 9)

for ($r = 2; $r -le 15; $r++) {
    foreach ($col in 5, 6) {
        $cell = $ws.Cells.Item($r, $col)
        $boolValue = $cell.Value2

        if ($col -eq 5 -and $quotedStyleRows -contains $r) {
            # Re-enter the value as literal text (leading apostrophe forces
            # a text cell instead of Excel's automatic TRUE/FALSE boolean).
            if ($boolValue) {
                $cell.Value = "'TRUE"
            } else {
                $cell.Value = "'FALSE"
            }
        } else {
            # Use a helper formula that yields the word as a text string,
            # then flatten it to a static value below.
            if ($boolValue) {
                $cell.Formula = "=""TRUE"""
            } else {
                $cell.Formula = "=""FALSE"""
            }
        }
    }
}

# Convert the helper formulas in E2:F15 into plain static text values
# (keeping their existing cell style) so no formula/quote-prefix is left
# behind.
$boolRange = $ws.Range("E2:F15")
$boolRange.Copy() | Out-Null
$boolRange.PasteSpecial(-4163) | Out-Null

# Update the sheet's active cell/selection.
$ws.Range("K10").Select() | Out-Null
